$d = $word.ActiveDocument

# Locate the "DEVLOG DAY 4" Report paragraph: "Hari ini tim kami " + "telah selesai ..."
# It currently consists of two runs; the target collapses them into a single run
# with the combined sentence.
$p14 = $d.Paragraphs.Item(14)
$mergeEnd = $p14.Range.End - 1
$mergeRng = $d.Range($p14.Range.Start, $mergeEnd)

$mergedRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr><w:t>Hari ini tim kami telah selesai membuat asset untuk kucing. Kami juga telah menguploadnya ke github dan mulai mencobanya untuk digerakan. Kami juga memasukan scene falling dimana scene ini akan ditrigger apabila anak kucing terjatuh ke air. Kami memastikan dan menguji coba asset yang kami masukan. Kedepannya kami mulai melengkapi beberapa bagian yang memerlukan asset dan menguji coba terus project yang kami buat.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mergeRng.InsertXML($mergedRunXml)

# Insert the new "DEVLOG DAY 5" heading, its "Report" sub-heading, and replace the
# text of the following (previously empty) paragraph with the new devlog entry.
$p15 = $d.Paragraphs.Item(15)
$insPoint = $d.Range($p15.Range.Start, $p15.Range.Start)

$newBlockXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-ID"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">DEVLOG DAY </w:t></w:r><w:r><w:rPr><w:lang w:val="en-ID"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:lang w:val="en-ID"/></w:rPr><w:br/><w:t>SALT STUDIO CODELABS</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr><w:t>Report</w:t></w:r></w:p>' +
'<w:p><w:r><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr><w:t xml:space="preserve">Hari ini tim kami </w:t></w:r><w:r><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr><w:t xml:space="preserve">telah menyelesaikan UI untuk game kami. Kebanyakan kami selesaikan dengan menggunakan figma. Kami </w:t></w:r><w:r><w:rPr><w:lang w:val="en-ID" w:eastAsia="en-ID"/></w:rPr><w:t>menyelesaikan UI untuk text, conveyor belt, cat paw dan juga asset untuk home. Asset ini masih kami uji coba untuk disatukan dengan menggunakan code. Selanjutnya kami akan mulai mempersatukan berbagai aspek agar dapat tercipta fisrt playable awal.</w:t></w:r></w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($newBlockXml)
